$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price and Volume columns remain text so values like "566.69" are not
# auto-converted to numbers by Excel when assigned via .Value
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "63.051.77"
$ws.Range("E2").Value = "  +0.25%  "

# Row 3
$ws.Range("D3").Value = "3.405.09"
$ws.Range("E3").Value = "  +1.52%  "

# Row 4
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").Value = "566.69"
$ws.Range("E5").Value = "  +1.63%  "

# Row 6
$ws.Range("D6").Value = "154.97"
$ws.Range("E6").Value = "  +1.04%  "

# Row 7
$ws.Range("E7").Value = "  +0.07%  "

# Row 8
$ws.Range("D8").Value = "3.407.84"
$ws.Range("E8").Value = "  +1.52%  "

# Row 9
$ws.Range("D9").Value = "0.541"
$ws.Range("E9").Value = "  +1.76%  "

# Row 10
$ws.Range("D10").Value = "7.39"
$ws.Range("E10").Value = "  -1.69%  "

# Row 11
$ws.Range("E11").Value = "  +2.48%  "

# Row 12
$ws.Range("D12").Value = "0.433"
$ws.Range("E12").Value = "  -1.59%  "

# Row 13
$ws.Range("D13").Value = "3.996.63"
$ws.Range("E13").Value = "  +1.64%  "

# Row 14
$ws.Range("E14").Value = "  -3.72%  "

# Row 15
$ws.Range("D15").Value = "0.0000188"
$ws.Range("E15").Value = "  +3.60%  "

# Row 16
$ws.Range("D16").Value = "26.87"
$ws.Range("E16").Value = "  -0.82%  "

# Row 17
$ws.Range("D17").Value = "63.209.14"
$ws.Range("E17").Value = "  +0.37%  "

# Row 18
$ws.Range("D18").Value = "3.433.10"
$ws.Range("E18").Value = "  +1.92%  "

# Row 19
$ws.Range("D19").Value = "6.24"
$ws.Range("E19").Value = "  -4.28%  "

# Row 20
$ws.Range("D20").Value = "14.03"
$ws.Range("E20").Value = "  +1.52%  "

# Row 21
$ws.Range("D21").Value = "382.77"
$ws.Range("E21").Value = "  -1.40%  "

# Row 22
$ws.Range("D22").Value = "8.09"
$ws.Range("E22").Value = "  -4.28%  "

# Row 23
$ws.Range("D23").Value = "0.996"
$ws.Range("E23").Value = "  -0.33%  "

# Row 24
$ws.Range("D24").Value = "71.45"
$ws.Range("E24").Value = "  +1.27%  "

# Row 25
$ws.Range("D25").Value = "0.530"
$ws.Range("E25").Value = "  -2.31%  "

# Row 26
$ws.Range("D26").Value = "0.0000117"
$ws.Range("E26").Value = "  +19.83%  "

# Row 27
$ws.Range("D27").Value = "9.39"
$ws.Range("E27").Value = "  +5.72%  "

# Row 28
$ws.Range("D28").Value = "0.175"
$ws.Range("E28").Value = "  -2.26%  "

# Row 29
$ws.Range("E29").Value = "  -0.11%  "

# Row 30
$ws.Range("D30").Value = "5.98"
$ws.Range("E30").Value = "  +5.79%  "

# Row 31
$ws.Range("D31").Value = "1.99"
$ws.Range("E31").Value = "  -0.08%  "

# Row 32
$ws.Range("D32").Value = "1.33"
$ws.Range("E32").Value = "  +1.78%  "

# Row 33
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "23.05"
$ws.Range("E33").Value = "  -0.22%  "

# Row 34
$ws.Range("B34").Value = "RenderToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D34").Value = "6.37"
$ws.Range("E34").Value = "  -3.77%  "

# Row 35
$ws.Range("E35").Value = "  +0.04%  "

# Row 36
$ws.Range("D36").Value = "6.76"
$ws.Range("E36").Value = "  +0.71%  "

# Row 37
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").Value = "157.88"
$ws.Range("E37").Value = "  -1.57%  "

# Row 38
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "1.45"
$ws.Range("E38").Value = "  -1.86%  "

# Row 39
$ws.Range("D39").Value = "0.0759"
$ws.Range("E39").Value = "  +2.25%  "

# Row 40
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "2.892.70"
$ws.Range("E40").Value = "  +2.31%  "

# Row 41
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "1.82"
$ws.Range("E41").Value = "  -3.75%  "

# Row 42
$ws.Range("D42").Value = "26.74"
$ws.Range("E42").Value = "  -1.36%  "

# Row 43
$ws.Range("D43").Value = "0.0316"
$ws.Range("E43").Value = "  +1.86%  "

# Row 44
$ws.Range("D44").Value = "4.33"
$ws.Range("E44").Value = "  +0.01%  "

# Row 45
$ws.Range("D45").Value = "0.756"
$ws.Range("E45").Value = "  +1.18%  "

# Row 46
$ws.Range("D46").Value = "41.04"
$ws.Range("E46").Value = "  +0.61%  "

# Row 47
$ws.Range("D47").Value = "23.38"
$ws.Range("E47").Value = "  +5.41%  "

# Row 48
$ws.Range("D48").Value = "1.07"
$ws.Range("E48").Value = "  +1.88%  "

# Row 49
$ws.Range("D49").Value = "2.15"
$ws.Range("E49").Value = "  +18.88%  "

# Row 50
$ws.Range("D50").Value = "6.38"
$ws.Range("E50").Value = "  +1.07%  "

# Row 51
$ws.Range("D51").Value = "0.831"
$ws.Range("E51").Value = "  +2.90%  "
